# Apply the "Trade #42 closed" update to the live trading results workbook.
#
# Summary:
#   - Summary sheet: update Total P&L %, Total Trades, Win Rate %
#   - Strategy Status sheet: update MarketMaking Trades, Win Rate %
#   - All Trades sheet: append new trade row (#42) as row 43
#   - MarketMaking sheet: append the same new trade row (#42) as row 43

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.1      # Total P&L %
$summary.Range("B6").Value = 42        # Total Trades
$summary.Range("B9").Value = 35.71     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 42         # MarketMaking Trades
$status.Range("G4").Value = 35.71      # MarketMaking Win Rate %

# ---------------------------------------------------------------------------
# 3) Helper: write the new trade-#42 row (row 43) into a given sheet
# ---------------------------------------------------------------------------
function Add-TradeRow42($ws) {
    $row = 43

    $ws.Cells.Item($row, 1).Value = 42          # A: Trade #

    # Use a leading apostrophe so Excel keeps these as literal text instead
    # of auto-converting the date-looking string into a date serial number.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"   # B: Date
    $ws.Cells.Item($row, 3).Value = "08:33:31"      # C: Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"  # D: Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"          # E: Side
    $ws.Cells.Item($row, 6).Value = 0.56            # F: Entry Price
    $ws.Cells.Item($row, 7).Value = 0.56            # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"        # H: Status
    $ws.Cells.Item($row, 9).Value = 0               # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0              # J: P&L $
    $ws.Cells.Item($row, 11).Value = 99.78          # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0              # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0              # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6            # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"   # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.13           # Q: Duration (min)
}

# ---------------------------------------------------------------------------
# 4) All Trades sheet - append row 43
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow42 $allTrades

# ---------------------------------------------------------------------------
# 5) MarketMaking sheet - append row 43
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow42 $marketMaking
